$d = $word.ActiveDocument

# 1. Trim the sentence about ASICs/market cap in the mining paragraph:
#    "... no ASICs will be created for quite some time, until Smartcash
#    reaches a considerable market cap." ->
#    "... no ASICs will be created for quite some time."
$findText    = "ASICs will be created for quite some time, until Smartcash reaches a considerable market cap."
$replaceText = "ASICs will be created for quite some time."

$editRange = $d.Content
$editRange.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2) | Out-Null

# Word automatically drops (or refreshes) a "_GoBack" bookmark at the site
# of the most recently edited text when the document is saved. Because the
# document already contains the "exchanges" bookmark further down, and
# bookmark ids are handed out in document order, this shifts "exchanges"
# from id 0 to id 1 - exactly the change seen in the target revision.
$d.Bookmarks.Add("_GoBack", $editRange) | Out-Null
